$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: log the 07.01.2024 session's worked minutes (GetDeckPlain work)
$ws.Range("B20").Value = 470

# Row 21: new work-log entry for 07.01.2024
# Set the work description first so its shared-string slot is allocated
# before the date string's slot (matches the original authoring order).
$ws.Range("C21").Value = "unit tests, protocol"

# Write the date as literal text (not an auto-converted date serial).
# Entering it via a formula that evaluates to a string, then collapsing
# the formula to its value with Paste Special (values only), keeps the
# cell a plain shared-string cell with the default (no) style - exactly
# like the other literal date cells in column A.
$ws.Range("A21").Formula = '="07.01.2024"'
$ws.Range("A21").Copy()
$ws.Range("A21").PasteSpecial(-4163)

# Reflect the new active selection left by the edits.
$ws.Range("B21").Select() | Out-Null
